# SurEau-Ecos calibration workbook update:
# - add an option for defoliation when PLC_TL > 10%: bump Depth_layer_1 (Feuil1!C2)
#   from 20 to 30 cm, which ripples through all the dependent formulas on both sheets.
# - leave the workbook positioned/scrolled/selected the way the author left it
#   (Feuil1 active, scrolled down, zoomed to 124%, selection on E27).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuil1")
$ws2 = $wb.Worksheets.Item("Feuil2")

# --- core data edit -------------------------------------------------------
$ws1.Range("C2").Value = 30

# --- window / view bookkeeping ---------------------------------------------
# Move the main Excel window close to where the author left it.
try {
    $excel.ActiveWindow.Left = 760
    $excel.ActiveWindow.Top = 500
} catch {}
try {
    $wb.Windows.Item(1).Left = 760
    $wb.Windows.Item(1).Top = 500
} catch {}

# Feuil1 becomes the active / selected tab (previously Feuil2 was active).
[void]$ws1.Activate()

# Scroll Feuil1 so row 18 / column A is the top-left visible cell, then zoom.
try {
    $excel.ActiveWindow.ScrollRow = 18
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
try {
    $excel.ActiveWindow.TopLeftCell = $ws1.Range("A18")
} catch {}

$excel.ActiveWindow.Zoom = 124

# Final selection on Feuil1 is E27.
[void]$ws1.Range("E27").Select()

# Feuil2's selection stays as it was (G7); only its tabSelected flag drops
# automatically now that Feuil1 is active.
